# Apply the "schedule.xlsx" update described by the commit:
#  - add two new columns (Vendor_class, max_products)
#  - refresh several Last_update / Next_update timestamps
#  - add a brand-new "Synology" vendor row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update existing date values
# ---------------------------------------------------------------------
$ws.Range("D3").Value = 45261

$ws.Range("C4").Value = 44899
$ws.Range("D4").Value = 44899

$ws.Range("C5").Value = 44899
$ws.Range("D5").Value = 44899

# ---------------------------------------------------------------------
# 2. Add the two new header cells (E1/F1), copying the header's
#    formatting (bold, bordered, centered) from an existing header cell
# ---------------------------------------------------------------------
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E1").Value = "Vendor_class"
$ws.Range("F1").Value = "max_products"

# ---------------------------------------------------------------------
# 3. Populate the new columns for the existing rows
# ---------------------------------------------------------------------
# Row 2 (Siemens) / Row 3 (Asus) -> new cells stay blank, nothing to set

# Row 4 (Schneider)
$ws.Range("E4").Value = "SchneiderElectricScraper"
$ws.Range("F4").Value = 10

# Row 5 (AVM)
$ws.Range("E5").Value = "AVMScraper"

# ---------------------------------------------------------------------
# 4. Add the brand-new row 6 (Synology)
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "Synology"
$ws.Range("B6").Value = 0

# Copy the date formatting from row 4 (C4 -> yyyy-mm-dd, D4 -> yyyy-mm-dd hh:mm:ss)
$ws.Range("C4:D4").Copy()
$ws.Range("C6:D6").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C6").Value = 44899
$ws.Range("D6").Value = 44899

$ws.Range("E6").Value = "Synology_scraper"

# ---------------------------------------------------------------------
# 5. Try to materialize the remaining "blank" string cells that the
#    source workbook stores explicitly (E2,F2,E3,F3,F5,F6). Some engines
#    drop truly-empty cells on save, so this is best effort.
# ---------------------------------------------------------------------
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("F6").Value = ""

Write-Host "Final UsedRange: $($ws.UsedRange.Address())"
